$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text number format on the blocks we are about to rewrite so that
# numeric/percentage-looking strings (e.g. "329.05", "0.63%") are stored as literal text
# instead of being auto-converted to numbers by Excel. Style is reset back to Normal
# afterwards so the cells end up unstyled again, exactly like the source file.
$block1 = $ws.Range("B2:E26")
$block1.NumberFormat = "@"
$block2 = $ws.Range("B38:E50")
$block2.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '329.05'
$ws.Range("E2").Value = '0.63%'

# Row 3
$ws.Range("D3").Value = '44.33'
$ws.Range("E3").Value = '0.71%'

# Row 4
$ws.Range("D4").Value = '5.598'
$ws.Range("E4").Value = '3.27%'

# Row 5
$ws.Range("D5").Value = '0.08089'
$ws.Range("E5").Value = '0.03%'

# Row 6
$ws.Range("D6").Value = '2.020'
$ws.Range("E6").Value = '6.43%'

# Row 7
$ws.Range("D7").Value = '4.315'
$ws.Range("E7").Value = '0.34%'

# Row 8
$ws.Range("D8").Value = '0.9536'
$ws.Range("E8").Value = '1.30%'

# Row 9
$ws.Range("E9").Value = '-5.13%'

# Row 10
$ws.Range("D10").Value = '0.1174'
$ws.Range("E10").Value = '-3.01%'

# Row 11
$ws.Range("D11").Value = '0.1858'
$ws.Range("E11").Value = '-1.78%'

# Row 12
$ws.Range("B12").Value = 'MCDex'
$ws.Range("C12").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D12").Value = '10.23'
$ws.Range("E12").Value = '20.52%'

# Row 13
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").Value = '0.09859'
$ws.Range("E13").Value = '3.32%'

# Row 14
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = '0.04597'
$ws.Range("E14").Value = '10.70%'

# Row 15
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = '0.1069'
$ws.Range("E15").Value = '-0.25%'

# Row 16
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = '0.001284'
$ws.Range("E16").Value = '-0.17%'

# Row 17
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = '0.04214'
$ws.Range("E17").Value = '-3.41%'

# Row 18
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").Value = '0.005930'
$ws.Range("E18").Value = '-2.14%'

# Row 19
$ws.Range("B19").Value = 'HotbitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D19").Value = '0.004322'
$ws.Range("E19").Value = '0.57%'

# Row 20
$ws.Range("B20").Value = 'LEO'
$ws.Range("C20").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D20").Value = '3.370'
$ws.Range("E20").Value = '-5.76%'

# Row 21
$ws.Range("B21").Value = 'BitpandaEcosystemToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D21").Value = '0.3475'
$ws.Range("E21").Value = '-0.67%'

# Row 22
$ws.Range("B22").Value = 'ProBitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D22").Value = '0.1410'
$ws.Range("E22").Value = '4.07%'

# Row 23
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").Value = '0.2505'
$ws.Range("E23").Value = '-3.81%'

# Row 24
$ws.Range("B24").Value = 'BitKan'
$ws.Range("C24").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D24").Value = '0.001249'
$ws.Range("E24").Value = '0.93%'

# Row 25
$ws.Range("D25").Value = '0.0001191'
$ws.Range("E25").Value = '-3.84%'

# Row 26
$ws.Range("E26").Value = '-0.91%'

# Row 38
$ws.Range("D38").Value = '0.02669'
$ws.Range("E38").Value = '1.13%'

# Row 39
$ws.Range("D39").Value = '0.05566'
$ws.Range("E39").Value = '2.42%'

# Row 40
$ws.Range("D40").Value = '0.007485'
$ws.Range("E40").Value = '-3.36%'

# Row 41
$ws.Range("E41").Value = '1.26%'

# Row 42
$ws.Range("D42").Value = '0.008074'
$ws.Range("E42").Value = '-17.22%'

# Row 43
$ws.Range("D43").Value = '0.002017'
$ws.Range("E43").Value = '-4.95%'

# Row 44
$ws.Range("D44").Value = '0.008405'
$ws.Range("E44").Value = '-15.06%'

# Row 45
$ws.Range("D45").Value = '0.00007232'
$ws.Range("E45").Value = '-1.30%'

# Row 46
$ws.Range("D46").Value = '0.00000000751'
$ws.Range("E46").Value = '-0.59%'

# Row 47
$ws.Range("D47").Value = '0.004196'
$ws.Range("E47").Value = '18.01%'

# Row 48
$ws.Range("D48").Value = '0.002269'
$ws.Range("E48").Value = '-0.69%'

# Row 49
$ws.Range("D49").Value = '0.00002102'
$ws.Range("E49").Value = '-0.59%'

# Row 50
$ws.Range("D50").Value = '0.0002001'
$ws.Range("E50").Value = '-0.59%'

# Restore default (unstyled) formatting on the two blocks
$block1.Style = "Normal"
$block2.Style = "Normal"